# Insert a new weekly data row at row 226 (pushes the existing rows 226..333
# down to 227..334, preserving their contents) and populate the new row with
# the latest "Choclo" (Choclero, Primera) price observation for the Maule
# region market.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(226).Insert()

$ws.Cells.Item(226, 1).Value = 5
$ws.Cells.Item(226, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(226, 3).Value = 'Maule'
$ws.Cells.Item(226, 4).Value = 45016
$ws.Cells.Item(226, 5).Value = 7
$ws.Cells.Item(226, 6).Value = 100112024
$ws.Cells.Item(226, 7).Value = 'Choclo'
$ws.Cells.Item(226, 8).Value = 'Choclero'
$ws.Cells.Item(226, 9).Value = 'Primera'
$ws.Cells.Item(226, 10).Value = 20000
$ws.Cells.Item(226, 11).Value = 400
$ws.Cells.Item(226, 12).Value = 400
$ws.Cells.Item(226, 13).Value = 400
$ws.Cells.Item(226, 14).Value = '$/unidad'
$ws.Cells.Item(226, 15).Value = 'Región del Maule'
$ws.Cells.Item(226, 16).Value = 400
$ws.Cells.Item(226, 17).Value = 1
$ws.Cells.Item(226, 18).Value = 'Hortaliza'
